# Automatic update from scheduled task (2025-07-26)
# - corrects the timestamp in A18 (tiny floating point refinement)
# - appends the new row 19 captured by the next scheduled run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: refine the stored timestamp value
$ws.Range("A18").Value = 45864.91690162037

# Row 19: new reading appended by the scheduled task
$ws.Range("A19").Value = 45864.95860632969
$ws.Range("A19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B19").Value = 2025
$ws.Range("C19").Value = 30
$ws.Range("D19").Value = 13.67
$ws.Range("E19").Value = 89.66
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 1.26
$ws.Range("H19").Value = "E"
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = "23:00:23"
